$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 54: fill OHLCV + SDs for SPY (cols C-H) ---
$ws.Range("C54").Value = 587.80999999999995
$ws.Range("D54").Value = 588.98
$ws.Range("E54").Value = 585.53
$ws.Range("F54").Value = 587.59
$ws.Range("G54").Value = 65952780
$ws.Range("H54").Value = 0.34561075894624332

# --- Row 54: HU54 (VT) value update 0 -> 585 ---
$ws.Range("HU54").Value = 585

# --- Row 55: new day of data (cols A, B, I..II) ---
$ws.Range("A55").Value = 45791
$ws.Range("B55").Value = 45792
$ws.Range("I55").Value = 0.1862
$ws.Range("J55").Value = 3.9
$ws.Range("K55").Value = 600
$ws.Range("L55").Value = 209518200
$ws.Range("M55").Value = 11786
$ws.Range("N55").Value = 2040
$ws.Range("O55").Value = 13826
$ws.Range("P55").Value = 0.080166234747757123
$ws.Range("Q55").Value = 0.040507848968097906
$ws.Range("R55").Value = 45793
$ws.Range("S55").Value = 0.14787501685283216
$ws.Range("T55").Value = 45807
$ws.Range("U55").Value = 0.10756304231141998
$ws.Range("V55").Value = 45828
$ws.Range("W55").Value = 0.27267755253631942
$ws.Range("X55").Value = 18.333333333333332
$ws.Range("Y55").Value = 590
$ws.Range("Z55").Value = 174819360
$ws.Range("AA55").Value = -442
$ws.Range("AB55").Value = 6270
$ws.Range("AC55").Value = 6712
$ws.Range("AD55").Value = 0.066889701478022728
$ws.Range("AE55").Value = 0.13343112277722022
$ws.Range("AF55").Value = 45792
$ws.Range("AG55").Value = 0.13343112277722022
$ws.Range("AH55").Value = 45793
$ws.Range("AI55").Value = 0.32630793381198692
$ws.Range("AJ55").Value = 45807
$ws.Range("AK55").Value = 0.11081088036233562
$ws.Range("AL55").Value = 6.333333333333333
$ws.Range("AM55").Value = 595
$ws.Range("AN55").Value = 104530195
$ws.Range("AO55").Value = 5665
$ws.Range("AP55").Value = 1117
$ws.Range("AQ55").Value = 6782
$ws.Range("AR55").Value = 0.039995533326454828
$ws.Range("AS55").Value = 0.096761167913777368
$ws.Range("AT55").Value = 45792
$ws.Range("AU55").Value = 0.096761167913777368
$ws.Range("AV55").Value = 45793
$ws.Range("AW55").Value = 0.28341899430200851
$ws.Range("AX55").Value = 45828
$ws.Range("AY55").Value = 0.1154204698973564
$ws.Range("AZ55").Value = 13.333333333333334
$ws.Range("BA55").Value = 605
$ws.Range("BB55").Value = 84411415
$ws.Range("BC55").Value = 12843
$ws.Range("BD55").Value = 123
$ws.Range("BE55").Value = 12966
$ws.Range("BF55").Value = 0.032297649131580675
$ws.Range("BG55").Value = 0
$ws.Range("BH55").Value = 45793
$ws.Range("BI55").Value = 0.072008739553887743
$ws.Range("BJ55").Value = 45807
$ws.Range("BK55").Value = 0.10297984224364592
$ws.Range("BL55").Value = 45828
$ws.Range("BM55").Value = 0.35388404043895272
$ws.Range("BN55").Value = 18.333333333333332
$ws.Range("BO55").Value = 610
$ws.Range("BP55").Value = 81032400
$ws.Range("BQ55").Value = 24649
$ws.Range("BR55").Value = 85
$ws.Range("BS55").Value = 24734
$ws.Range("BT55").Value = 0.03100476426665633
$ws.Range("BU55").Value = 0
$ws.Range("BV55").Value = 45828
$ws.Range("BW55").Value = 0.26094278275062505
$ws.Range("BX55").Value = 45856
$ws.Range("BY55").Value = 0.12148140471779992
$ws.Range("BZ55").Value = 45919
$ws.Range("CA55").Value = 0.13102193090405068
$ws.Range("CB55").Value = 76.666666666666671
$ws.Range("CC55").Value = 586
$ws.Range("CD55").Value = -140358134
$ws.Range("CE55").Value = 0.060742676398089585
$ws.Range("CF55").Value = 965
$ws.Range("CG55").Value = 4045
$ws.Range("CH55").Value = 5010
$ws.Range("CI55").Value = 0.56077382618691674
$ws.Range("CJ55").Value = 45792
$ws.Range("CK55").Value = 0.56077382618691674
$ws.Range("CL55").Value = 45793
$ws.Range("CM55").Value = 0.34993716235935984
$ws.Range("CN55").Value = 45796
$ws.Range("CO55").Value = 0.022588791670885625
$ws.Range("CP55").Value = 2.6666666666666665
$ws.Range("CQ55").Value = 587
$ws.Range("CR55").Value = -126199130
$ws.Range("CS55").Value = 0.054615095661719459
$ws.Range("CT55").Value = 1116
$ws.Range("CU55").Value = 5046
$ws.Range("CV55").Value = 6162
$ws.Range("CW55").Value = 0.30641386871640103
$ws.Range("CX55").Value = 45792
$ws.Range("CY55").Value = 0.30641386871640103
$ws.Range("CZ55").Value = 45793
$ws.Range("DA55").Value = 0.54866547203621618
$ws.Range("DB55").Value = 45796
$ws.Range("DC55").Value = 0.030009195510704519
$ws.Range("DD55").Value = 2.6666666666666665
$ws.Range("DE55").Value = 585
$ws.Range("DF55").Value = -109447065.00000001
$ws.Range("DG55").Value = 0.047365318008685389
$ws.Range("DH55").Value = 2289
$ws.Range("DI55").Value = 31210
$ws.Range("DJ55").Value = 33499
$ws.Range("DK55").Value = 0.26250847589443471
$ws.Range("DL55").Value = 45792
$ws.Range("DM55").Value = 0.26250847589443471
$ws.Range("DN55").Value = 45793
$ws.Range("DO55").Value = 0.36628015777923068
$ws.Range("DP55").Value = 45807
$ws.Range("DQ55").Value = 0.077825115867603489
$ws.Range("DR55").Value = 6.333333333333333
$ws.Range("DS55").Value = 550
$ws.Range("DT55").Value = -70473700
$ws.Range("DU55").Value = 0.030498846284719387
$ws.Range("DV55").Value = -889
$ws.Range("DW55").Value = 17152
$ws.Range("DX55").Value = 18041
$ws.Range("DY55").Value = 0
$ws.Range("DZ55").Value = 45807
$ws.Range("EA55").Value = 0.096201214223764098
$ws.Range("EB55").Value = 45828
$ws.Range("EC55").Value = 0.54428943129723706
$ws.Range("ED55").Value = 45856
$ws.Range("EE55").Value = 0.2219254119687771
$ws.Range("EF55").Value = 39.333333333333336
$ws.Range("EG55").Value = 575
$ws.Range("EH55").Value = -62744575
$ws.Range("EI55").Value = 0.027153919094996388
$ws.Range("EJ55").Value = -1800
$ws.Range("EK55").Value = 44338
$ws.Range("EL55").Value = 46138
$ws.Range("EM55").Value = 0.10283328144173183
$ws.Range("EN55").Value = 45793
$ws.Range("EO55").Value = 0.12381458215051007
$ws.Range("EP55").Value = 45807
$ws.Range("EQ55").Value = 0.21867731360231499
$ws.Range("ER55").Value = 45828
$ws.Range("ES55").Value = 0.20340286807933922
$ws.Range("ET55").Value = 18.333333333333332
$ws.Range("EU55").Value = 580
$ws.Range("EV55").Value = 419661900
$ws.Range("EW55").Value = -12600
$ws.Range("EX55").Value = 20146
$ws.Range("EY55").Value = 32746
$ws.Range("EZ55").Value = 0.085223563995324714
$ws.Range("FA55").Value = 180254140
$ws.Range("FB55").Value = 0.068969166886194505
$ws.Range("FC55").Value = 0.024711776384165156
$ws.Range("FD55").Value = 45793
$ws.Range("FE55").Value = 0.22900866521013055
$ws.Range("FF55").Value = 45828
$ws.Range("FG55").Value = 0.46001550921382445
$ws.Range("FH55").Value = 45919
$ws.Range("FI55").Value = 0.048059256780454529
$ws.Range("FJ55").Value = 55.666666666666664
$ws.Range("FK55").Value = -239407760
$ws.Range("FL55").Value = 0.10360830312029864
$ws.Range("FM55").Value = 0.092506274650412335
$ws.Range("FN55").Value = 45793
$ws.Range("FO55").Value = 0.33889895632455691
$ws.Range("FP55").Value = 45800
$ws.Range("FQ55").Value = 0.11910691616679427
$ws.Range("FR55").Value = 45807
$ws.Range("FS55").Value = 0.096634461639839908
$ws.Range("FT55").Value = 9
$ws.Range("FU55").Value = 585
$ws.Range("FV55").Value = 371962305
$ws.Range("FW55").Value = 2289
$ws.Range("FX55").Value = 31210
$ws.Range("FY55").Value = 33499
$ws.Range("FZ55").Value = 0.075536886488899724
$ws.Range("GA55").Value = 131257620
$ws.Range("GB55").Value = 0.0502220292907819
$ws.Range("GC55").Value = 0.090403437149020383
$ws.Range("GD55").Value = 45792
$ws.Range("GE55").Value = 0.090403437149020383
$ws.Range("GF55").Value = 45793
$ws.Range("GG55").Value = 0.36025885582871303
$ws.Range("GH55").Value = 45800
$ws.Range("GI55").Value = 0.10587773875528141
$ws.Range("GJ55").Value = 4
$ws.Range("GK55").Value = -240704685
$ws.Range("GL55").Value = 0.1041695723060773
$ws.Range("GM55").Value = 0.26250847589443471
$ws.Range("GN55").Value = 45792
$ws.Range("GO55").Value = 0.26250847589443471
$ws.Range("GP55").Value = 45793
$ws.Range("GQ55").Value = 0.36628015777923068
$ws.Range("GR55").Value = 45807
$ws.Range("GS55").Value = 0.077825115867603489
$ws.Range("GT55").Value = 6.333333333333333
$ws.Range("GU55").Value = 590
$ws.Range("GV55").Value = 283716840
$ws.Range("GW55").Value = -442
$ws.Range("GX55").Value = 6270
$ws.Range("GY55").Value = 6712
$ws.Range("GZ55").Value = 0.057616286516100942
$ws.Range("HA55").Value = 229268100
$ws.Range("HB55").Value = 0.087722977406126307
$ws.Range("HC55").Value = 0.13343112277722022
$ws.Range("HD55").Value = 45792
$ws.Range("HE55").Value = 0.13343112277722022
$ws.Range("HF55").Value = 45793
$ws.Range("HG55").Value = 0.32630793381198692
$ws.Range("HH55").Value = 45807
$ws.Range("HI55").Value = 0.11081088036233562
$ws.Range("HJ55").Value = 6.333333333333333
$ws.Range("HK55").Value = -54448740
$ws.Range("HL55").Value = 0.023563737275844064
$ws.Range("HM55").Value = 0
$ws.Range("HN55").Value = 45807
$ws.Range("HO55").Value = 0.094954814381379621
$ws.Range("HP55").Value = 45828
$ws.Range("HQ55").Value = 0.16351342565502894
$ws.Range("HR55").Value = 45856
$ws.Range("HS55").Value = 0.20582753613765903
$ws.Range("HT55").Value = 39.333333333333336
$ws.Range("HU55").Value = 588
$ws.Range("HV55").Value = 40420
$ws.Range("HW55").Value = 201298
$ws.Range("HX55").Value = 2613546721.5
$ws.Range("HY55").Value = -2310700521
$ws.Range("HZ55").Value = 302846200.5
$ws.Range("IA55").Value = 1.1310625058278594
$ws.Range("IB55").Value = 4924247242.5
$ws.Range("IC55").Value = 0.13870530324006167
$ws.Range("ID55").Value = 45792
$ws.Range("IE55").Value = 0.13870530324006167
$ws.Range("IF55").Value = 45793
$ws.Range("IG55").Value = 0.25731944023117359
$ws.Range("IH55").Value = 45828
$ws.Range("II55").Value = 0.15101016914464974

# --- Update the selection in the bottom-right (frozen) pane ---
$ws.Range("E61").Select()
